# Apply the "LAYOUT IS BETTER FOR SQL DATABASE" edit:
#  1. Fill every blank/empty data cell (rows 2-37, i.e. excluding the header
#     row) with the literal text "NULL" so the sheet is friendlier to load
#     into a SQL database.
#  2. Update the sheet view: zoom to 60% and change the active selection in
#     the bottom-right (frozen) pane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$lastCol = $usedRange.Column + $usedRange.Columns.Count - 1

# Row 1 is the header row and must stay untouched; data starts at row 2.
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($null -eq $val -or $val -eq "") {
            $cell.Value = "NULL"
        }
    }
}

# Adjust the sheet view: zoom level and the selected range in the scrolled
# (bottom-right) pane.
$win = $excel.ActiveWindow
$win.Zoom = 60
$ws.Range("M19:V19").Select()
